$d = $word.ActiveDocument

# Insert a "." run at the very start of the first paragraph (before "Here, we have raw dataset on")
$p1 = $d.Paragraphs(1)
$startRange = $p1.Range.Duplicate
$startRange.Collapse(1)  # wdCollapseStart
$startRange.InsertBefore(".")

# Insert a "." run at the very end of the first paragraph's text (after " down below.")
$p1 = $d.Paragraphs(1)
$endRange = $p1.Range.Duplicate
$endRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark (wdCharacter)
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertAfter(".")
